# Updates cryptos list prices / 1h volume percentages (and fixes the
# ARBITRUM / TrustWalletToken row ordering) to match the latest scrape.
#
# Note: several Price values in column D look like plain decimal numbers
# (e.g. "209.68"), and Excel's COM automation auto-converts such strings
# to numeric values when assigned directly via .Value. To preserve them
# as text (matching the original inline-string cell content), those
# cells are explicitly formatted as Text ("@") before assignment, then
# the format is reset back to the workbook's default "Normal" style so
# no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.834.71'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.584.78'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('E4').Value = '  +0.07%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '209.68'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.36%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.478'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -3.82%  '
$ws.Range('E8').Value = '  -0.90%  '
$ws.Range('E9').Value = '  -0.41%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '18.02'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -2.18%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0792'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.805.74'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '1.585.87'
$ws.Range('E13').Value = '  -2.60%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '4.02'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').Value = '25.825.53'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').Value = '0.0₃0723'
$ws.Range('E17').Value = '  -2.08%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '59.70'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -3.25%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '191.61'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('E24').Value = '  -0.96%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '142.05'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('E26').Value = '  +0.01%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '1.71'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.35%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.10'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.92%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.45'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -2.93%  '
$ws.Range('E30').Value = '  -5.73%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.0471'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('E32').Value = '  -0.40%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.02'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.53%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D36').Value = '1.100.08'
$ws.Range('E36').Value = '  -2.55%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  -2.14%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.502'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.820'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +8.10%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.775'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -8.26%  '
$ws.Range('E43').Value = '  +1.48%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '93.81'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -4.01%  '
$ws.Range('D45').Value = '1.719.32'
$ws.Range('E45').Value = '  -2.17%  '
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('E47').Value = '  -0.62%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '53.24'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('E49').Value = '  -1.63%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.407'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('E51').Value = '  +0.01%  '
